# Update database and change read_price algorithm
# - Insert 5 new quarterly columns (D:H) before the existing data, shifting
#   the previously-reported quarters from D:H -> I:M.
# - Fill the newly inserted columns with the earlier-period data (period
#   labels, publish dates, and the financial figures for each line item).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Insert five new columns before column D. Excel shifts the existing
#    D:H content (and its number formats / styles) to I:M automatically.
# ---------------------------------------------------------------------------
$ws.Range("D:H").Insert()

# Re-apply column widths for the newly inserted D:H block, matching the
# same convention used throughout the sheet: every column is 29 characters
# wide except the one showing the "فصل چهارم" (4th-quarter / longest
# ordinal) label, which is widened to 31 so the text isn't clipped. For the
# new block that longest label lands in column F.
$ws.Range("D1").ColumnWidth = 28.17
$ws.Range("E1").ColumnWidth = 28.17
$ws.Range("F1").ColumnWidth = 30.17
$ws.Range("G1").ColumnWidth = 28.17
$ws.Range("H1").ColumnWidth = 28.17

# ---------------------------------------------------------------------------
# 2) Header row 8 - new quarter period labels for the newly inserted columns
# ---------------------------------------------------------------------------
$ws.Range("D8").Value = "فصل دوم منتهی به 1399/06"
$ws.Range("E8").Value = "فصل سوم منتهی به 1399/09"
$ws.Range("F8").Value = "فصل چهارم منتهی به 1399/12"
$ws.Range("G8").Value = "فصل اول منتهی به 1400/03"
$ws.Range("H8").Value = "فصل دوم منتهی به 1400/06"

# ---------------------------------------------------------------------------
# 3) Header row 9 - publish dates for the newly inserted columns
# ---------------------------------------------------------------------------
$ws.Range("D9").Value = "1400-09-30 (4)"
$ws.Range("E9").Value = "1400-10-30 (2)"
$ws.Range("F9").Value = "1401-04-15 (8)"
$ws.Range("G9").Value = "1401-04-29 (2)"
$ws.Range("H9").Value = "1401-09-14 (4)"

# ---------------------------------------------------------------------------
# 4) Financial data rows 11-27 for the newly inserted columns D:H
# ---------------------------------------------------------------------------
$data = @{
    11 = @(27803150, 39823676, 50066714, 49950179, 26869807)
    12 = @(-16518619, -26338713, -41112264, -40185447, -20614487)
    13 = @(11284531, 13484963, 8954450, 9764732, 6255320)
    14 = @(-403532, -560697, -869903, -747548, -524652)
    15 = @(0, 0, 0, 0, 0)
    16 = @(-168929, 1365819, 14898, -424567, -628134)
    17 = @(10712070, 14290085, 8099445, 8592617, 5102534)
    18 = @(-344507, -1311191, -680741, -1022293, -609200)
    19 = @(184738, -6512, 4565849, -349894, -51091)
    20 = @(10552301, 12972382, 11984553, 7220430, 4442243)
    21 = @(-1072342, 0, -2593718, 0, -662248)
    22 = @(9479959, 12972382, 9390835, 7220430, 3779995)
    23 = @(0, 0, 0, 0, 0)
    24 = @(9479959, 12972382, 9390835, 7220430, 3779995)
    25 = @(632, 865, 626, 481, 118)
    26 = @(15000000, 15000000, 15000000, 15000000, 32000000)
    27 = @(296, 405, 293, 226, 118)
}

$cols = @("D", "E", "F", "G", "H")

foreach ($row in $data.Keys) {
    $values = $data[$row]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range($cols[$i] + $row).Value = $values[$i]
    }
}
